# Rename the three embedded logo inline pictures:
#   - BTec_Logo-Orange (first-page header):            image2.jpg -> image1.jpg
#   - Pearson logo (default footer):                   image1.png -> image2.png
#   - Pearson logo (first-page footer):                image1.png -> image2.png
#
# Each picture's range is selected before the rename - going through
# $word.Selection.InlineShapes avoids a stale-handle resolution issue that
# otherwise hits InlineShapes fetched straight off a HeaderFooter.Range
# (most noticeable for footers).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page header holds the BTec logo.
$btecHeader = $sec.Headers.Item(2)
$btecInline = $btecHeader.Range.InlineShapes.Item(1)
$btecInline.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

# Default (primary) footer holds one Pearson logo instance.
$defaultFooter = $sec.Footers.Item(1)
$pearsonDefaultInline = $defaultFooter.Range.InlineShapes.Item(1)
$pearsonDefaultInline.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# First-page footer holds the other Pearson logo instance.
$firstFooter = $sec.Footers.Item(2)
$pearsonFirstInline = $firstFooter.Range.InlineShapes.Item(1)
$pearsonFirstInline.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"
